$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "797/DR TADLA"
$ws.Range("B2").Value = "Direction régionale"
$ws.Range("C2").Value = "ad646456"
$ws.Range("D2").Value = "Mimo crimo"
$ws.Range("E2").Value = "non"
$ws.Range("F2").Value = "mensuelle"
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = "--"
$ws.Range("I2").Value = 6000
$ws.Range("J2").Value = "--"
$ws.Range("K2").Value = 600
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 5400

# Row 3
$ws.Range("A3").Value = "010/DR010/AV"
$ws.Range("B3").Value = "Direction régionale"
$ws.Range("C3").Value = "aa654556"
$ws.Range("D3").Value = "Ali Ali"
$ws.Range("E3").Value = "non"
$ws.Range("F3").Value = "mensuelle"
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = "--"
$ws.Range("I3").Value = 2000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 2000

# Row 4 (new)
$ws.Range("A4").Value = "001/DR TANGER/AV"
$ws.Range("B4").Value = "Direction régionale"
$ws.Range("C4").Value = "ada666"
$ws.Range("D4").Value = "Mohamed berrada"
$ws.Range("E4").Value = "non"
$ws.Range("F4").Value = "mensuelle"
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = "--"
$ws.Range("I4").Value = 5000
$ws.Range("J4").Value = "--"
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 4500

# Row 5 (new)
$ws.Range("A5").Value = "000/DR DEV"
$ws.Range("B5").Value = "Direction régionale"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "117946464"
$ws.Range("D5").Value = "IBM"
$ws.Range("E5").Value = "oui"
$ws.Range("F5").Value = "mensuelle"
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 30000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 20000
$ws.Range("M5").Value = 30000

# Row 6 (new)
$ws.Range("A6").Value = "000/DR DEV"
$ws.Range("B6").Value = "Direction régionale"
$ws.Range("C6").Value = "BJ179134"
$ws.Range("D6").Value = "Ahmed Tawfiq"
$ws.Range("E6").Value = "non"
$ws.Range("F6").Value = "mensuelle"
$ws.Range("G6").Value = 15
$ws.Range("H6").Value = 30000
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 1500
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 20000
$ws.Range("M6").Value = 28500

# Row 7 (was row 4, total row - shifted down and values changed)
$ws.Range("A7").Value = " "
$ws.Range("B7").Value = " "
$ws.Range("C7").Value = " "
$ws.Range("D7").Value = " "
$ws.Range("E7").Value = " "
$ws.Range("F7").Value = " "
$ws.Range("G7").Value = " "
$ws.Range("H7").Value = 60000
$ws.Range("I7").Value = 13000
$ws.Range("J7").Value = 1500
$ws.Range("K7").Value = 1100
$ws.Range("L7").Value = 40000
$ws.Range("M7").Value = 70400
